$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 453) holds the "Förändrad" date stamp, which was
# bumped from 2023-09-13 (serial 45182) to 2023-09-15 (serial 45184) for
# every data row in the sheet.
$ws.Range("C2:C453").Value = 45184
